$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the mailto: hyperlinks on the email column (C2/C3) while keeping
# their "Hyperlink" cell style intact.
$ws.Hyperlinks.Delete() | Out-Null

# Drop the second registrant row entirely (bulk regis -> single fresh row).
$ws.Rows.Item(3).Delete() | Out-Null

# Reset the remaining registrant's data (reset password & new applicant).
$ws.Range("I2").Value = "documents/blank.pdf"
$ws.Range("A2").Value = "'198409142003121002"
$ws.Range("B2").Value = "ASEP RIYANTO, S.STP."
$ws.Range("C2").Value = "asriy.oppof9@gmail.com"
$ws.Range("D2").Value = "''082316840508"
$ws.Range("E2").Value = "Pemerintah Provinsi Jawa Barat"
$ws.Range("F2").Value = "Analis SDM Aparatur"
$ws.Range("G2").Value = "Madya"
$ws.Range("H2").Value = "documents/jabarkolektif9.jpeg"

# Re-fit the columns to the new content widths.
$ws.Cells.EntireColumn.AutoFit() | Out-Null

$ws.Range("H23").Select() | Out-Null
